$wb = $excel.ActiveWorkbook

# Add the two new sheets at the end, in order: methods, architecture
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "methods"
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "architecture"

# ----- methods sheet -----
$ws4.Range("A1").Value = "/boats"
$ws4.Range("B1").Value = "/users"
$ws4.Range("C1").Value = "/loads"

$ws4.Range("A3").Value = "GET"
$ws4.Range("B3").Value = "GET"
$ws4.Range("C3").Value = "GET"

$ws4.Range("A5").Value = "/boats/<boat_id>"
$ws4.Range("B5").Value = "/users/<user_id>"
$ws4.Range("C5").Value = "/loads/<load_id>"

$ws4.Range("A7").Value = "GET"
$ws4.Range("B7").Value = "GET"
$ws4.Range("C7").Value = "GET"

$ws4.Range("A8").Value = "POST"
$ws4.Range("B8").Value = "POST"
$ws4.Range("C8").Value = "POST"

$ws4.Range("A9").Value = "PUT"
$ws4.Range("B9").Value = "PUT"
$ws4.Range("C9").Value = "PUT"

$ws4.Range("A10").Value = "PATCH"
$ws4.Range("B10").Value = "PATCH"
$ws4.Range("C10").Value = "PATCH"

$ws4.Range("A11").Value = "DELETE"
$ws4.Range("B11").Value = "DELETE"
$ws4.Range("C11").Value = "DELETE"

# column widths for methods
$ws4.Columns.Item(1).ColumnWidth = 35
$ws4.Columns.Item(2).ColumnWidth = 44.5
$ws4.Columns.Item(3).ColumnWidth = 46

# font color black for column C rows 7-11
$ws4.Range("C7:C11").Font.Color = 0

$ws4.Range("B13").Select()

# ----- architecture sheet -----
$ws5.Range("A1").Value = "models"
$ws5.Range("B1").Value = "controllers"
$ws5.Range("C1").Value = "auth"

$ws5.Range("F29").Select()
